# Apply updated crypto price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.604.07"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.089.83"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'591.36"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'154.85"
$ws.Range("E6").Value = "  +6.80%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  +2.61%  "

$ws.Range("D9").Value = "3.083.13"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("E10").Value = "  -1.77%  "

$ws.Range("D11").Value = "'5.84"

$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").Value = "'37.44"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "'0.0000240"
$ws.Range("E14").Value = "  -2.00%  "

$ws.Range("D15").Value = "3.601.44"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Value = "63.579.22"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "3.083.95"
$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").Value = "'476.54"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("D21").Value = "'14.60"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("D22").Value = "'0.714"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'2.41"
$ws.Range("E24").Value = "  +4.01%  "

$ws.Range("D25").Value = "'12.90"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("D26").Value = "'81.04"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  +2.05%  "

$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "'2.68"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "'2.18"
$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "  +4.32%  "

$ws.Range("D34").Value = "'27.18"
$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").Value = "0.0₃0848"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  -1.50%  "

$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  +5.22%  "

$ws.Range("D38").Value = "'6.07"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").Value = "'2.22"
$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("D40").Value = "'9.28"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").Value = "'50.68"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").Value = "'444.00"
$ws.Range("E42").Value = "  -1.95%  "

$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").Value = "'0.0362"
$ws.Range("E44").Value = "  -2.55%  "

$ws.Range("D45").Value = "'40.01"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("D47").Value = "2.800.87"
$ws.Range("E47").Value = "  -3.92%  "

$ws.Range("D48").Value = "'131.27"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").Value = "'25.73"
$ws.Range("E49").Value = "  +6.87%  "

$ws.Range("E51").Value = "  +1.20%  "

